# TeacherDetails.xlsx update
#
# 1. Delete "Sheet1" (the tiny Branch/Andheri/123 helper sheet) entirely.
# 2. Rename the remaining "Sheet2" (the real teacher-data sheet) to "Sheet1",
#    taking over the tab-selected / active state.
# 3. Merge the old Branch value pair ("Andheri" / "123") into a single cell
#    value "123,Mumbai Branch" in the cell that referenced "Andheri" (S2).

$wb = $excel.ActiveWorkbook

$oldSheet1 = $wb.Worksheets.Item("Sheet1")
$oldSheet1.Delete()

$dataSheet = $wb.Worksheets.Item("Sheet2")
$dataSheet.Name = "Sheet1"

$dataSheet.Range("S2").Value = "123,Mumbai Branch"

$dataSheet.Range("A1").Select()
